# Apply the "made changes to some graphs" edit:
# - Column D: rows where value is 184.02 get corrected to 184
# - Row 29: restore the "normal" precision values (matching the other rows)
# - Move the active selection to B30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column D whose value should be corrected from 184.02 to 184
$rowsToFix = @(7, 8, 10, 12, 13, 15, 17, 19, 20, 22, 23, 25, 27, 28)
foreach ($r in $rowsToFix) {
    $ws.Range("D$r").Value = 184
}

# Row 29 gets the same precision as the rest of the table
$ws.Range("A29").Value = 2754.42
$ws.Range("B29").Value = 243.88
$ws.Range("C29").Value = 183.18
$ws.Range("D29").Value = 184
$ws.Range("E29").Value = 182.35

# Update the active selection/cursor position
$ws.Range("B30").Select()
